# Refresh the cryptocurrency price table (columns B-E) to the latest scrape.
# Column A (rank index) is untouched. A few coins swapped adjacent rank rows
# (rows 36/37, 39/40, 46/47), so their full row (Coin/Link/Price/Volume) moves.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") is stored as plain text in the source data (values like
# "60.027.15" or "0.420" are not valid numbers / would lose a trailing zero),
# so format it as Text before writing to stop Excel from auto-converting it.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "60.027.15"
$ws.Cells.Item(2, 5).Value = "  +3.56%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.421.53"
$ws.Cells.Item(3, 5).Value = "  +3.04%  "

$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "552.81"
$ws.Cells.Item(5, 5).Value = "  +2.09%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "138.34"
$ws.Cells.Item(6, 5).Value = "  +2.67%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.578"
$ws.Cells.Item(8, 5).Value = "  +1.65%  "

$ws.Cells.Item(9, 5).Value = "  +3.16%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "5.75"
$ws.Cells.Item(10, 5).Value = "  +3.96%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.358"
$ws.Cells.Item(11, 5).Value = "  +0.38%  "

$ws.Cells.Item(12, 5).Value = "  -2.01%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "24.90"
$ws.Cells.Item(13, 5).Value = "  +4.65%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "2.854.32"
$ws.Cells.Item(14, 5).Value = "  +3.04%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "59.983.52"
$ws.Cells.Item(15, 5).Value = "  +3.62%  "

$ws.Cells.Item(16, 5).Value = "  +2.99%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.416.80"
$ws.Cells.Item(17, 5).Value = "  +1.91%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "11.35"
$ws.Cells.Item(18, 5).Value = "  +5.94%  "

$ws.Cells.Item(19, 5).Value = "  +2.09%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "331.68"
$ws.Cells.Item(20, 5).Value = "  +0.54%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.79"
$ws.Cells.Item(21, 5).Value = "  +0.82%  "

$ws.Cells.Item(22, 5).Value = "  -0.02%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "65.06"
$ws.Cells.Item(23, 5).Value = "  +3.63%  "

$ws.Cells.Item(24, 5).Value = "  +3.42%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "8.59"
$ws.Cells.Item(25, 5).Value = "  +2.79%  "

$ws.Cells.Item(26, 5).Value = "  +0.08%  "

$ws.Cells.Item(27, 5).Value = "  -0.72%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.0₃0786"
$ws.Cells.Item(28, 5).Value = "  +6.65%  "

$ws.Cells.Item(29, 5).Value = "  +0.92%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "6.28"
$ws.Cells.Item(30, 5).Value = "  +2.40%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "168.67"
$ws.Cells.Item(31, 5).Value = "  -0.68%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.05"
$ws.Cells.Item(32, 5).Value = "  +2.38%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "18.68"
$ws.Cells.Item(33, 5).Value = "  +1.63%  "

$ws.Cells.Item(34, 5).Value = "  -0.02%  "

$ws.Cells.Item(35, 5).Value = "  +4.89%  "

$ws.Cells.Item(36, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.00"
$ws.Cells.Item(36, 5).Value = "  +0.18%  "

$ws.Cells.Item(37, 2).Value = "NEARProtocol"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.21"
$ws.Cells.Item(37, 5).Value = "  +0.18%  "

$ws.Cells.Item(38, 5).Value = "  -0.02%  "

$ws.Cells.Item(39, 2).Value = "Bittensor"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "322.55"
$ws.Cells.Item(39, 5).Value = "  +11.66%  "

$ws.Cells.Item(40, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.420"
$ws.Cells.Item(40, 5).Value = "  +11.12%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "39.45"
$ws.Cells.Item(41, 5).Value = "  +0.93%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "3.69"
$ws.Cells.Item(42, 5).Value = "  +1.11%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "139.88"
$ws.Cells.Item(43, 5).Value = "  -1.93%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.0961"

$ws.Cells.Item(45, 5).Value = "  +2.07%  "

$ws.Cells.Item(46, 2).Value = "Polygon"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.415"
$ws.Cells.Item(46, 5).Value = "  +8.84%  "

$ws.Cells.Item(47, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "19.46"
$ws.Cells.Item(47, 5).Value = "  +1.76%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.574"
$ws.Cells.Item(48, 5).Value = "  +1.33%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.0226"
$ws.Cells.Item(49, 5).Value = "  +1.77%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "17.76"
$ws.Cells.Item(50, 5).Value = "  +1.42%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "11.06"

